$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "blood_pressure" column (H) entirely; this shifts the
# "hospital_expire_flag" column (I) left into H, matching the diff.
$ws.Columns.Item(8).Delete()

# Move selection to match the post-edit sheet view (F8).
$ws.Range("F8").Select()
